$wb = $excel.ActiveWorkbook

$wsValid = $wb.Worksheets.Item("ValidLogin")
$wsInvalid = $wb.Worksheets.Item("InvalidLogin")

# Remove the now-invalid extra "valid login" test rows (3,4,5) from the
# InvalidLogin sheet, shifting the remaining data up.
$wsInvalid.Rows.Item(3).Delete()
$wsInvalid.Rows.Item(3).Delete()
$wsInvalid.Rows.Item(3).Delete()

# Update the remaining credentials row with the new invalid-login test data.
$wsInvalid.Range("A2").Value = "mobiliya1234@gmail.com"
$wsInvalid.Range("B2").Value = "mobiliya123"

# Select B3 on the InvalidLogin sheet, which also makes it the active tab
# (matching the "waits" / active-sheet change from the commit).
$wsInvalid.Range("B3").Select()

# Give the columns a bit more breathing room on both sheets.
$wsValid.Columns.Item(1).ColumnWidth = 27.605442176870767
$wsValid.Columns.Item(2).ColumnWidth = 24.365646258503368
$wsValid.Columns.Item(3).ColumnWidth = 13.748299319727867

$wsInvalid.Columns.Item(1).ColumnWidth = 24.54421768707487
$wsInvalid.Columns.Item(2).ColumnWidth = 13.926870748299367
$wsInvalid.Columns.Item(3).ColumnWidth = 12.125850340136067
